$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("CustomerEmail" shifts right, etc.)
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Channel"

# Match column A's "best fit" width (Excel carries the left neighbour's
# width into a freshly inserted column).
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Place the selection on the newly inserted header cell, matching the
# post-edit state captured by the author's save.
$ws.Range("B1").Select()
